$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.989.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.758.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5200"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2708"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06220"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.758.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07017"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.502"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.979.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006706"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.980.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.082"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.390"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.215"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.480"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.819"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08381"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.697"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.417"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04409"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.721"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01568"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.958"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3887"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7486"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.926"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05491"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1117"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.095"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "

